# Applies the coin-tracker refresh described in the commit:
#   "Updated cryptos list on Thu Oct 19 05:05:42 UTC 2023 with GitHub Actions"
# Price (col D) and Volume(1h) (col E) cells hold plain text, not numbers
# (values like "209.62" or "1.48" must stay text, matching the sheet's
# existing inlineStr cells) - so any Price cell whose new value still looks
# like a number is forced to Text format before the write to stop Excel
# from re-interpreting it as a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '28.306.45'
$ws.Range("E2").Value = '  -1.42%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.551.25'
$ws.Range("E3").Value = '  -1.32%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.76'
$ws.Range("E5").Value = '  -1.60%  '

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.483'
$ws.Range("E6").Value = '  -1.77%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.01%  '

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.242'
$ws.Range("E9").Value = '  -2.16%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -1.42%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.772.79'
$ws.Range("E12").Value = '  -1.38%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.544.66'
$ws.Range("E13").Value = '  -1.63%  '

# Row 14: WrappedBTC
$ws.Range("D14").Value = '28.291.83'
$ws.Range("E14").Value = '  -1.39%  '

# Row 15: Polkadot
$ws.Range("E15").Value = '  -1.70%  '

# Row 16: Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.509'
$ws.Range("E16").Value = '  -2.61%  '

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.54'
$ws.Range("E17").Value = '  -3.05%  '

# Row 18: BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.49'
$ws.Range("E18").Value = '  -1.88%  '

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.32'
$ws.Range("E19").Value = '  -0.95%  '

# Row 20: ShibaInu
$ws.Range("E20").Value = '  -2.91%  '

# Row 21: Dai
$ws.Range("E21").Value = '  -0.02%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  +0.14%  '

# Row 23: Avalanche
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.86'
$ws.Range("E23").Value = '  -3.26%  '

# Row 24: Toncoin
$ws.Range("E24").Value = '  -6.01%  '

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.11'
$ws.Range("E25").Value = '  -1.00%  '

# Row 26: EthereumClassic
$ws.Range("E26").Value = '  -1.96%  '

# Row 27: Stellar
$ws.Range("E27").Value = '  -1.14%  '

# Row 28: BinanceUSD
$ws.Range("E28").Value = '  +0.10%  '

# Row 29: Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.24'
$ws.Range("E29").Value = '  -3.40%  '

# Row 30: Hedera
$ws.Range("E30").Value = '  -4.12%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -4.75%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  -1.65%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = '  -2.35%  '

# Row 34: Maker
$ws.Range("D34").Value = '1.385.33'
$ws.Range("E34").Value = '  -0.52%  '

# Row 35: TrustWalletToken
$ws.Range("E35").Value = '  +0.58%  '

# Row 36: LidoDAOToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  -3.70%  '

# Row 37: HuobiToken
$ws.Range("E37").Value = '  -1.24%  '

# Row 38: MXToken
$ws.Range("E38").Value = '  -0.63%  '

# Row 39: VeChain
$ws.Range("E39").Value = '  -3.31%  '

# Row 40: ImmutableX
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("E40").Value = '  +0.91%  '

# Row 41: RenderToken
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.509'
$ws.Range("E41").Value = '  -3.78%  '

# Row 42: PaxDollar
$ws.Range("E42").Value = '  +0.08%  '

# Row 43: ARBITRUM
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.776'
$ws.Range("E43").Value = '  -2.03%  '

# Row 44: Kaspa
$ws.Range("E44").Value = '  -0.84%  '

# Row 45: FraxShare
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.39'
$ws.Range("E45").Value = '  -2.29%  '

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.82'
$ws.Range("E46").Value = '  -2.30%  '

# Row 47: WEMIXToken
$ws.Range("E47").Value = '  -6.08%  '

# Row 48: RocketPoolETH
$ws.Range("D48").Value = '1.686.23'
$ws.Range("E48").Value = '  -1.41%  '

# Row 49: Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.61'
$ws.Range("E49").Value = '  -1.22%  '

# Row 50: BitcoinSV
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '41.94'
$ws.Range("E50").Value = '  +2.52%  '

# Row 51: BabyDogeCoin
$ws.Range("E51").Value = '  +3.22%  '
